$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New data row to append (next day's report after the last existing row, 06/14/2020 -> 06/15/2020)
$newRow = 95
$values = @(43997, 1023, 320, 464, 258, 54)

# Carry the formatting down from the row above (date number format in col A,
# centered alignment in B:F) before writing the new values
$ws.Range("A94:F94").Copy()
$ws.Range("A95:F95").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 1).Value = $values[0]
$ws.Cells.Item($newRow, 2).Value = $values[1]
$ws.Cells.Item($newRow, 3).Value = $values[2]
$ws.Cells.Item($newRow, 4).Value = $values[3]
$ws.Cells.Item($newRow, 5).Value = $values[4]
$ws.Cells.Item($newRow, 6).Value = $values[5]

# Resize the table (ListObject) so the autofilter/table range covers the new row
$table = $ws.ListObjects.Item("Condicion_Pacientes")
$table.Resize($ws.Range("A1:F95"))

# Update the view: scroll so the new last row is in frame and select it,
# mirroring where the workbook was left after entering the new data
$excel.ActiveWindow.ScrollRow = 84
$ws.Range("F95").Select()
